$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for 8e59d3a0-...md file is now ready for handoff ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# Latest Handoff Date got refreshed when the handoff report was regenerated
$overview.Range("D2").Value = "2016-27-18 03:27:15"
$overview.Range("D3").Value = "2016-27-18 03:27:15"

# --- zh-cn detail sheet ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("E2").Value = "2016-03-18 03:27:07"
$zh.Range("E3").Value = "2016-03-18 03:27:07"

# --- de-de detail sheet ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = "Ready for handoff"
$de.Range("E2").Value = "2016-03-18 03:27:15"
$de.Range("E3").Value = "2016-03-18 03:27:15"
